$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    2,
    -0.3036867844311208,
    0.226016869225491,
    0.1547881991105919,
    -0.06892048399615523,
    0.1588074336107037,
    -0.4459064887304521,
    -0.2533299822585643,
    -0.3869446842432583,
    0.3689886109514408,
    -0.2375676495630193,
    -0.06821648309498324,
    -0.009372778487127869,
    0.09835750693104157,
    0.0367282424715893,
    0.4193125732117537,
    0.6483189581467216,
    0.05431395812032949,
    0.4194873667791191,
    0.3597478966727853,
    0.5589686541212187,
    0.4224199210089787,
    -0.02969368961751345,
    4.2964596053585,
    0.4655444954022276,
    0.3692330055649125,
    0.255308055126365,
    0.985866018235656,
    5.244289870569874,
    0.9343558192232184,
    -0.2778607738411689,
    0.7490073304897067,
    0.885522784570417,
    -0.7030225518903331,
    0.7996699859181595,
    0.7480588837663376,
    0.7194746032728944,
    0.7058849505409454,
    0.5694553353364057,
    0.7387958335202249,
    0.5510444178321311,
    0.6647413438913339,
    0.6990383841299922,
    0.6635160940718288,
    0.6310103914282831,
    -1.279920598506624,
    -0.9941628244536368,
    -0.8876851908265895,
    -0.6530785597146361,
    -0.05576828822634386,
    -0.8808236328114679,
    -0.8808236328114679,
    -1.114756251233469,
    -0.1941042360003489,
    -1.013853905987674,
    -0.8976590279862675,
    -0.9287350184708151,
    -1.123567071138395,
    -0.8504870520082407,
    -0.4903224098973188,
    0.3718422085082914,
    -1.202259018807272,
    -0.7343364021664399,
    -0.8974519606912046,
    -0.106159846411747,
    -0.7851092787761388,
    -0.7594586154370843
)

# Header cell C1: copy style from B1 (bold header style) and set value 2
$ws.Range("C1").Value = $values[0]
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Fill in C2:C67 with the predicted 't+3' values
for ($i = 1; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
